$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight metrics ("Natmi following Dr Hou advice").
# Each entry is Cell,Value as produced by the re-run of the NATMI pipeline.
$updates = @(
    ,@("E2", 3)
    ,@("G2", 2.264243666666667)
    ,@("H2", 6.792731)
    ,@("I2", 0.4114976873616865)
    ,@("J2", 0.4114976873616865)
    ,@("K2", 3)
    ,@("M2", 7.602918333333334)
    ,@("N2", 22.808755)
    ,@("O2", 0.4364663819587059)
    ,@("P2", 0.4364663819587059)
    ,@("Q2", 17.21485968443389)
    ,@("R2", 154.933737159905)
    ,@("S2", 0.17960490678713)
    ,@("T2", 0.17960490678713)
    ,@("E3", 3)
    ,@("G3", 2.264243666666667)
    ,@("H3", 6.792731)
    ,@("I3", 0.4114976873616865)
    ,@("J3", 0.4114976873616865)
    ,@("K3", 3)
    ,@("M3", 3.132521333333333)
    ,@("N3", 9.397564000000001)
    ,@("O3", 0.1798309797402525)
    ,@("P3", 0.1798309797402525)
    ,@("Q3", 7.092791589698223)
    ,@("R3", 63.835124307284)
    ,@("S3", 0.07400003227910022)
    ,@("T3", 0.07400003227910022)
    ,@("E4", 3)
    ,@("G4", 2.264243666666667)
    ,@("H4", 6.792731)
    ,@("I4", 0.4114976873616865)
    ,@("J4", 0.4114976873616865)
    ,@("K4", 3)
    ,@("M4", 0.6975466666666668)
    ,@("N4", 2.09264)
    ,@("O4", 0.04004457979149087)
    ,@("P4", 0.04004457979149087)
    ,@("Q4", 1.579415622204445)
    ,@("R4", 14.21474059984)
    ,@("S4", 0.01647825197556902)
    ,@("T4", 0.01647825197556902)
    ,@("E5", 3)
    ,@("G5", 2.264243666666667)
    ,@("H5", 6.792731)
    ,@("I5", 0.4114976873616865)
    ,@("J5", 0.4114976873616865)
    ,@("K5", 3)
    ,@("M5", 5.986266666666666)
    ,@("N5", 17.9588)
    ,@("O5", 0.3436580585095507)
    ,@("P5", 0.3436580585095507)
    ,@("Q5", 13.55436638697778)
    ,@("R5", 121.9892974828)
    ,@("S5", 0.1414144963198873)
    ,@("T5", 0.1414144963198873)
    ,@("E6", 3)
    ,@("G6", 2.089228666666667)
    ,@("H6", 6.267686)
    ,@("I6", 0.3796909216792509)
    ,@("J6", 0.3796909216792509)
    ,@("K6", 3)
    ,@("M6", 7.602918333333334)
    ,@("N6", 22.808755)
    ,@("O6", 0.4364663819587059)
    ,@("P6", 0.4364663819587059)
    ,@("Q6", 15.88423493232556)
    ,@("R6", 142.95811439093)
    ,@("S6", 0.165722322847909)
    ,@("T6", 0.165722322847909)
    ,@("E7", 3)
    ,@("G7", 2.089228666666667)
    ,@("H7", 6.267686)
    ,@("I7", 0.3796909216792509)
    ,@("J7", 0.3796909216792509)
    ,@("K7", 3)
    ,@("M7", 3.132521333333333)
    ,@("N7", 9.397564000000001)
    ,@("O7", 0.1798309797402525)
    ,@("P7", 0.1798309797402525)
    ,@("Q7", 6.544553368544889)
    ,@("R7", 58.90098031690401)
    ,@("S7", 0.06828019044405918)
    ,@("T7", 0.06828019044405918)
    ,@("E8", 3)
    ,@("G8", 2.089228666666667)
    ,@("H8", 6.267686)
    ,@("I8", 0.3796909216792509)
    ,@("J8", 0.3796909216792509)
    ,@("K8", 3)
    ,@("M8", 0.6975466666666668)
    ,@("N8", 2.09264)
    ,@("O8", 0.04004457979149087)
    ,@("P8", 0.04004457979149087)
    ,@("Q8", 1.457334492337778)
    ,@("R8", 13.11601043104)
    ,@("S8", 0.01520456340928947)
    ,@("T8", 0.01520456340928947)
    ,@("E9", 3)
    ,@("G9", 2.089228666666667)
    ,@("H9", 6.267686)
    ,@("I9", 0.3796909216792509)
    ,@("J9", 0.3796909216792509)
    ,@("K9", 3)
    ,@("M9", 5.986266666666666)
    ,@("N9", 17.9588)
    ,@("O9", 0.3436580585095507)
    ,@("P9", 0.3436580585095507)
    ,@("Q9", 12.50667992631111)
    ,@("R9", 112.5601193368)
    ,@("S9", 0.1304838449779932)
    ,@("T9", 0.1304838449779932)
    ,@("E10", 3)
    ,@("G10", 0.1905406666666667)
    ,@("H10", 0.571622)
    ,@("I10", 0.03462835949856721)
    ,@("J10", 0.03462835949856721)
    ,@("K10", 3)
    ,@("M10", 7.602918333333334)
    ,@("N10", 22.808755)
    ,@("O10", 0.4364663819587059)
    ,@("P10", 0.4364663819587059)
    ,@("Q10", 1.448665127845556)
    ,@("R10", 13.03798615061)
    ,@("S10", 0.01511411478350502)
    ,@("T10", 0.01511411478350502)
    ,@("E11", 3)
    ,@("G11", 0.1905406666666667)
    ,@("H11", 0.571622)
    ,@("I11", 0.03462835949856721)
    ,@("J11", 0.03462835949856721)
    ,@("K11", 3)
    ,@("M11", 3.132521333333333)
    ,@("N11", 9.397564000000001)
    ,@("O11", 0.1798309797402525)
    ,@("P11", 0.1798309797402525)
    ,@("Q11", 0.5968727032008889)
    ,@("R11", 5.371854328808)
    ,@("S11", 0.006227251815425021)
    ,@("T11", 0.006227251815425021)
    ,@("E12", 3)
    ,@("G12", 0.1905406666666667)
    ,@("H12", 0.571622)
    ,@("I12", 0.03462835949856721)
    ,@("J12", 0.03462835949856721)
    ,@("K12", 3)
    ,@("M12", 0.6975466666666668)
    ,@("N12", 2.09264)
    ,@("O12", 0.04004457979149087)
    ,@("P12", 0.04004457979149087)
    ,@("Q12", 0.1329110068977778)
    ,@("R12", 1.19619906208)
    ,@("S12", 0.001386678104988806)
    ,@("T12", 0.001386678104988806)
    ,@("E13", 3)
    ,@("G13", 0.1905406666666667)
    ,@("H13", 0.571622)
    ,@("I13", 0.03462835949856721)
    ,@("J13", 0.03462835949856721)
    ,@("K13", 3)
    ,@("M13", 5.986266666666666)
    ,@("N13", 17.9588)
    ,@("O13", 0.3436580585095507)
    ,@("P13", 0.3436580585095507)
    ,@("Q13", 1.140627241511111)
    ,@("R13", 10.2656451736)
    ,@("S13", 0.01190031479464837)
    ,@("T13", 0.01190031479464837)
    ,@("E14", 3)
    ,@("G14", 0.9584326666666668)
    ,@("H14", 2.875298)
    ,@("I14", 0.1741830314604954)
    ,@("J14", 0.1741830314604954)
    ,@("K14", 3)
    ,@("M14", 7.602918333333334)
    ,@("N14", 22.808755)
    ,@("O14", 0.4364663819587059)
    ,@("P14", 0.4364663819587059)
    ,@("Q14", 7.286885292665557)
    ,@("R14", 65.58196763399002)
    ,@("S14", 0.07602503754016189)
    ,@("T14", 0.07602503754016188)
    ,@("E15", 3)
    ,@("G15", 0.9584326666666668)
    ,@("H15", 2.875298)
    ,@("I15", 0.1741830314604954)
    ,@("J15", 0.1741830314604954)
    ,@("K15", 3)
    ,@("M15", 3.132521333333333)
    ,@("N15", 9.397564000000001)
    ,@("O15", 0.1798309797402525)
    ,@("P15", 0.1798309797402525)
    ,@("Q15", 3.002310774896889)
    ,@("R15", 27.02079697407201)
    ,@("S15", 0.03132350520166812)
    ,@("T15", 0.03132350520166812)
    ,@("E16", 3)
    ,@("G16", 0.9584326666666668)
    ,@("H16", 2.875298)
    ,@("I16", 0.1741830314604954)
    ,@("J16", 0.1741830314604954)
    ,@("K16", 3)
    ,@("M16", 0.6975466666666668)
    ,@("N16", 2.09264)
    ,@("O16", 0.04004457979149087)
    ,@("P16", 0.04004457979149087)
    ,@("Q16", 0.6685515118577779)
    ,@("R16", 6.016963606720002)
    ,@("S16", 0.006975086301643574)
    ,@("T16", 0.006975086301643574)
    ,@("E17", 3)
    ,@("G17", 0.9584326666666668)
    ,@("H17", 2.875298)
    ,@("I17", 0.1741830314604954)
    ,@("J17", 0.1741830314604954)
    ,@("K17", 3)
    ,@("M17", 5.986266666666666)
    ,@("N17", 17.9588)
    ,@("O17", 0.3436580585095507)
    ,@("P17", 0.3436580585095507)
    ,@("Q17", 5.737433524711112)
    ,@("R17", 51.6369017224)
    ,@("S17", 0.05985940241702185)
    ,@("T17", 0.05985940241702185)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = [double]$u[1]
}

